# Fruta / hortaliza, semanal
# Insert a new weekly record at row 306 (pushing the existing rows 306-330
# down to 307-331) and populate the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 306:330 down to 307:331, carrying formatting/styles with them.
$ws.Rows.Item(306).Insert()

# Populate the newly inserted row 306 with the new weekly record.
$ws.Cells.Item(306, 1).Value  = 7
$ws.Cells.Item(306, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(306, 3).Value  = "Ñuble"
$ws.Cells.Item(306, 4).Value  = 44826
$ws.Cells.Item(306, 5).Value  = 16
$ws.Cells.Item(306, 6).Value  = 100114013
$ws.Cells.Item(306, 7).Value  = "Zanahoria"
$ws.Cells.Item(306, 8).Value  = "Sin especificar"
$ws.Cells.Item(306, 9).Value  = "Primera"
$ws.Cells.Item(306, 10).Value = 120
$ws.Cells.Item(306, 11).Value = 11000
$ws.Cells.Item(306, 12).Value = 12000
$ws.Cells.Item(306, 13).Value = 11500
$ws.Cells.Item(306, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(306, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(306, 16).Value = 575
$ws.Cells.Item(306, 17).Value = 20
$ws.Cells.Item(306, 18).Value = "Hortaliza"
